{"js": "const results = context.document.body.search(\"North Yorkshire, UK\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Yorkshire, UK\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"North Yorkshire, UK\", $true, $false, $false, $false, $false, $true, 1, $false, \"Yorkshire, UK\", 2)\n"}
